# "Label BOM items better."
# Rename a handful of Package/Description entries in the BOM sheet so the
# labelling is clearer (resistor wattage, capacitor lead-spacing, and the
# ceramic/electrolytic capacitor descriptions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resistor package label: "R-W4" -> "R-1/4W"
$ws.Range("C2").Value = "R-1/4W"
$ws.Range("C3").Value = "R-1/4W"
$ws.Range("C4").Value = "R-1/4W"
$ws.Range("C5").Value = "R-1/4W"

# Ceramic capacitor package label: "C-5mm" -> "C-P5mm"
$ws.Range("C6").Value = "C-P5mm"

# Electrolytic capacitor package label: "E2.5-6.3" -> "E-P2.5mm 6.3x11.5mm"
$ws.Range("C7").Value = "E-P2.5mm 6.3x11.5mm"

# Ceramic capacitor description: "Capacitor Ceramic THT" -> "Ceramic Capacitor THT"
$ws.Range("E6").Value = "Ceramic Capacitor THT"

# Electrolytic capacitor description: "Capacitor Polarized THT" -> "Electrolytic Capacitor THT"
$ws.Range("E7").Value = "Electrolytic Capacitor THT"

# Update the active selection left by the author after making the edits
$ws.Range("E6:E7").Select()
